$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 38), columns G:J mirror the existing header row 22 (C:F)
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# Row 39: label + averages of column F (Flow_Lac), mirroring row 23 pattern (B column)
$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = '=AVERAGE(F$1:F$3)'
$ws.Range("H39").Formula = '=AVERAGE(F$4:F$6)'
$ws.Range("I39").Formula = '=AVERAGE(F$9:F$11)'
$ws.Range("J39").Formula = '=AVERAGE(F$13:F$16)'

# Row 40: standard error of the mean for column F, mirroring row 24 pattern
$ws.Range("G40").Formula = '=STDEV(F$1:F$3)/SQRT(COUNT(F$1:F$3))'
$ws.Range("H40").Formula = '=STDEV(F$4:F$6)/SQRT(COUNT(F$4:F$6))'
$ws.Range("I40").Formula = '=STDEV(F$9:F$11)/SQRT(COUNT(F$9:F$11))'
$ws.Range("J40").Formula = '=STDEV(F$13:F$16)/SQRT(COUNT(F$13:F$16))'

# Update selection/view to match the new working area (mirrors the author
# scrolling down to, and selecting, the newly added block)
$ws.Range("F38:J40").Select()
